$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.876.22'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.10%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.723.76'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.76%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.76'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.523'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.34%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.77'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +11.32%  '

$ws.Range("E9").Value = '  +3.61%  '

$ws.Range("E10").Value = '  +1.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0901'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.966.26'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.713.26'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("E14").Value = '  +3.19%  '

$ws.Range("E15").Value = '  +5.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.94'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.852.06'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '241.51'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.53%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0754'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.33%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.93'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.17%  '

$ws.Range("E21").Value = '  -0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.63'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.20%  '

$ws.Range("E23").Value = '  +4.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.47'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.51'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.74%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.73'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.10%  '

$ws.Range("E28").Value = '  +1.41%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0509'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.11%  '

$ws.Range("E31").Value = '  +1.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.46'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.93%  '

$ws.Range("E33").Value = '  +3.79%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.473.92'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.18%  '

$ws.Range("E35").Value = '  -2.04%  '

$ws.Range("E36").Value = '  +4.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.612'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.69%  '

$ws.Range("E38").Value = '  +0.76%  '

$ws.Range("E39").Value = '  -0.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.12'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '71.40'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.51%  '

$ws.Range("E42").Value = '  +6.27%  '

$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.870.69'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.84%  '

$ws.Range("E45").Value = '  +0.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.791'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.71'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '91.58'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.94%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0111'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.78%  '

$ws.Range("E50").Value = '  +2.36%  '

$ws.Range("E51").Value = '  +0.38%  '
